# project_finances.xlsx — "added purchasing list xlsx"
#
# The exchange-rate assumption in C1 moved from 3 to 3.1, so the
# "Price MYR" column (I) formulas were repointed from the hard-coded
# literal 3 to the assumption cell $C$1, and the column's running total
# (I19) was switched from a continuation of the per-row shared formula to
# an explicit SUM so it keeps working independently of the per-row range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assume 1 SGD = 3.1 MYR (was 3)
$ws.Range("C1").Value = 3.1

# Price MYR = Price SGD (H) * exchange rate (now $C$1 instead of literal 3)
$ws.Range("I9").Formula = '=H9*$C$1'
$ws.Range("I10:I18").Formula = '=H10*$C$1'

# Total row: sum the column instead of extending the per-row formula
$ws.Range("I19").Formula = "=SUM(I9:I18)"

# Restore the view/selection state recorded in the saved workbook
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("L19").Select()

$wb.Application.Calculate()
